$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 1.826271637843651
$ws.Cells.Item(3, 3).Value = 1.906592789843651
$ws.Cells.Item(4, 3).Value = 1.910950668843651
$ws.Cells.Item(5, 3).Value = 2.452713729843651
$ws.Cells.Item(6, 3).Value = 1.806423198843651
$ws.Cells.Item(7, 2).Value = 1.454170430843651
$ws.Cells.Item(8, 2).Value = 1.527312638843651
$ws.Cells.Item(9, 2).Value = -0.2024072381563491
$ws.Cells.Item(46, 4).Value = 0.4006729710360878
$ws.Cells.Item(47, 4).Value = 0.3913241300360878
$ws.Cells.Item(48, 4).Value = 0.07541333303608777
$ws.Cells.Item(49, 4).Value = 0.1149542130360878
$ws.Cells.Item(50, 3).Value = 0.7012562470360878
$ws.Cells.Item(51, 3).Value = 1.029332457036088
$ws.Cells.Item(52, 3).Value = 1.485428559036088
$ws.Cells.Item(53, 3).Value = 1.831943276036088
$ws.Cells.Item(53, 4).Value = 0.55793444658209
$ws.Cells.Item(54, 3).Value = 1.696990870036088
$ws.Cells.Item(54, 4).Value = 0.44417418258209
$ws.Cells.Item(55, 3).Value = 0.1642582080360878
$ws.Cells.Item(55, 4).Value = -0.54231891241791
$ws.Cells.Item(56, 2).Value = 0.1757658360360879
$ws.Cells.Item(56, 3).Value = -0.63179670641791
$ws.Cells.Item(57, 2).Value = 0.1903804690360879
$ws.Cells.Item(57, 3).Value = 0.24128467758209
$ws.Cells.Item(58, 3).Value = 0.44433217858209
$ws.Cells.Item(59, 3).Value = 0.6062157845820899
$ws.Cells.Item(59, 4).Value = -0.8596988570317647
$ws.Cells.Item(60, 3).Value = 1.06983264558209
$ws.Cells.Item(60, 4).Value = -0.9590976180317647
$ws.Cells.Item(61, 3).Value = 0.7401392785820899
$ws.Cells.Item(61, 4).Value = -1.043273225031765
$ws.Cells.Item(62, 2).Value = 0.65042024058209
$ws.Cells.Item(62, 3).Value = -0.5948165710317647
$ws.Cells.Item(63, 2).Value = 0.66923223158209
$ws.Cells.Item(63, 3).Value = -1.459370673031765
$ws.Cells.Item(64, 3).Value = -1.543924081031765
$ws.Cells.Item(65, 3).Value = -1.197149625031765
$ws.Cells.Item(66, 3).Value = -0.8452485510317647
$ws.Cells.Item(67, 3).Value = -1.137602302031765
$ws.Cells.Item(68, 2).Value = -1.037127551031765
$ws.Cells.Item(69, 2).Value = -1.085847230031765
$ws.Cells.Item(71, 4).Value = 0.7977233685636995
$ws.Cells.Item(72, 4).Value = 1.031128409563699
$ws.Cells.Item(73, 4).Value = 0.9920349685636995
$ws.Cells.Item(74, 4).Value = 1.0168253855637
$ws.Cells.Item(75, 3).Value = 1.053976176563699
$ws.Cells.Item(76, 3).Value = 1.093632598563699
$ws.Cells.Item(77, 3).Value = 1.124545940563699
$ws.Cells.Item(78, 3).Value = 1.0707846905637
$ws.Cells.Item(78, 4).Value = 0.8985202055291455
$ws.Cells.Item(79, 3).Value = 0.9333254135636995
$ws.Cells.Item(79, 4).Value = 0.8577491755291454
$ws.Cells.Item(80, 3).Value = 0.9579434665636994
$ws.Cells.Item(80, 4).Value = 0.8987580675291454
$ws.Cells.Item(81, 2).Value = 1.0751219075637
$ws.Cells.Item(81, 3).Value = 1.258920946529146
$ws.Cells.Item(82, 2).Value = 1.1210615525637
$ws.Cells.Item(82, 3).Value = 1.493423065529146
$ws.Cells.Item(83, 3).Value = 1.027544699529146
$ws.Cells.Item(84, 3).Value = 0.5469191505291455
$ws.Cells.Item(84, 4).Value = -0.9817231827224345
$ws.Cells.Item(85, 3).Value = 0.6114425455291455
$ws.Cells.Item(85, 4).Value = -0.9854454907224345
$ws.Cells.Item(86, 3).Value = 0.3885492905291454
$ws.Cells.Item(86, 4).Value = -0.8873425837224346
$ws.Cells.Item(87, 3).Value = -0.1343109554708546
$ws.Cells.Item(87, 4).Value = -1.085683986722434
$ws.Cells.Item(88, 2).Value = 0.4230266025291455
$ws.Cells.Item(88, 3).Value = 0.1374695422775655
$ws.Cells.Item(89, 2).Value = 0.4182579295291455
$ws.Cells.Item(89, 3).Value = -0.2331624027224345
$ws.Cells.Item(90, 3).Value = -0.03326235772243452
$ws.Cells.Item(91, 3).Value = -0.2778009377224345
$ws.Cells.Item(91, 4).Value = -1.165608299123972
$ws.Cells.Item(92, 3).Value = -0.6454920347224345
$ws.Cells.Item(92, 4).Value = -1.024285167559779
$ws.Cells.Item(93, 3).Value = -0.7103774527224346
$ws.Cells.Item(93, 4).Value = -1.000794259827642
$ws.Cells.Item(94, 2).Value = -0.5011116027224345
$ws.Cells.Item(94, 3).Value = -0.6420067559859775
$ws.Cells.Item(95, 2).Value = -0.5544081717224345
$ws.Cells.Item(95, 3).Value = -0.7670271480143395
$ws.Cells.Item(96, 3).Value = -0.9225503716806988
$ws.Cells.Item(97, 3).Value = -0.8491698657853378
$ws.Cells.Item(97, 4).Value = 0.8526545954887239
$ws.Cells.Item(98, 3).Value = -1.191515643655161
$ws.Cells.Item(98, 4).Value = 1.016949629488724
$ws.Cells.Item(99, 3).Value = -1.430441087857995
$ws.Cells.Item(99, 4).Value = 0.8815376954887238
$ws.Cells.Item(100, 2).Value = -1.217909980957737
$ws.Cells.Item(100, 3).Value = 0.8721537754887239
$ws.Cells.Item(101, 2).Value = -1.26539928353432
$ws.Cells.Item(101, 3).Value = 0.5557457034887239
$ws.Cells.Item(102, 3).Value = 0.6390669014887239
$ws.Cells.Item(103, 3).Value = 0.5246197534887239
$ws.Cells.Item(103, 4).Value = 0.3082097950934801
$ws.Cells.Item(104, 3).Value = 0.2428672904887239
$ws.Cells.Item(104, 4).Value = 0.3048056840934801
$ws.Cells.Item(105, 3).Value = 0.5514316234887239
$ws.Cells.Item(105, 4).Value = 0.4343567680934801
$ws.Cells.Item(106, 2).Value = 0.3148484034887239
$ws.Cells.Item(106, 3).Value = 0.9884114590934802
$ws.Cells.Item(107, 2).Value = 0.3295457764887239
$ws.Cells.Item(107, 3).Value = 1.11054283609348
$ws.Cells.Item(108, 3).Value = 0.6209472400934801
$ws.Cells.Item(109, 3).Value = 0.9007734940934802
$ws.Cells.Item(109, 4).Value = 0.7305722247131936
$ws.Cells.Item(110, 3).Value = 0.5169714940934801
$ws.Cells.Item(110, 4).Value = 0.8421929667131937
$ws.Cells.Item(111, 3).Value = 0.2115338810934801
$ws.Cells.Item(111, 4).Value = 0.4788276257131937
$ws.Cells.Item(112, 2).Value = 0.0506311690934801
$ws.Cells.Item(112, 3).Value = -0.1957891702868064
$ws.Cells.Item(113, 2).Value = -0.02970344090651991
$ws.Cells.Item(113, 3).Value = -0.5800827092868064
$ws.Cells.Item(114, 3).Value = -0.2004644112868063
$ws.Cells.Item(115, 3).Value = -0.1799067152868063
$ws.Cells.Item(116, 3).Value = -0.4927747992868063
$ws.Cells.Item(117, 3).Value = -0.3419294162868063
$ws.Cells.Item(119, 2).Value = -0.3059158432868063
$ws.Cells.Item(120, 2).Value = -0.2670276532868063
$ws.Cells.Item(122, 4).Value = -0.4083682634916527
$ws.Cells.Item(123, 4).Value = -0.5334118554916527
$ws.Cells.Item(124, 4).Value = -0.4617206544916527
$ws.Cells.Item(125, 4).Value = -0.7651034874916527
$ws.Cells.Item(126, 3).Value = -0.2379622684916527
$ws.Cells.Item(127, 3).Value = -0.4135428994916527
$ws.Cells.Item(128, 3).Value = -0.4776197014916527
$ws.Cells.Item(129, 3).Value = -0.2036544774916527
$ws.Cells.Item(129, 4).Value = 0.5987332491758083
$ws.Cells.Item(130, 3).Value = -0.1998718194916527
$ws.Cells.Item(130, 4).Value = 0.6101592851758083
$ws.Cells.Item(131, 3).Value = -0.1834133014916527
$ws.Cells.Item(131, 4).Value = 0.5298148361758083
$ws.Cells.Item(132, 3).Value = 0.03791119950834732
$ws.Cells.Item(132, 4).Value = 0.5915373291758083
$ws.Cells.Item(133, 2).Value = -0.05564365149165268
$ws.Cells.Item(133, 3).Value = 0.4130432101758083
$ws.Cells.Item(134, 2).Value = 0.03353349250834725
$ws.Cells.Item(134, 3).Value = -0.3288647778241918
$ws.Cells.Item(135, 3).Value = 0.04982731217580827
$ws.Cells.Item(136, 3).Value = 0.3369026561758083
$ws.Cells.Item(136, 4).Value = 0.7254492243564907
$ws.Cells.Item(137, 3).Value = 0.3645469811758083
$ws.Cells.Item(137, 4).Value = 0.7215746373564907
$ws.Cells.Item(138, 3).Value = 0.1236970551758083
$ws.Cells.Item(138, 4).Value = 0.5311946523564907
$ws.Cells.Item(139, 3).Value = 0.2095627611758083
$ws.Cells.Item(139, 4).Value = 0.5539812373564907
$ws.Cells.Item(140, 2).Value = -0.08176241982419175
$ws.Cells.Item(140, 3).Value = 0.1753415943564907
$ws.Cells.Item(141, 2).Value = -0.1256759188241917
$ws.Cells.Item(141, 3).Value = 0.2651053283564908
$ws.Cells.Item(142, 3).Value = 0.08763596535649075
$ws.Cells.Item(143, 3).Value = 0.1003532183564907
$ws.Cells.Item(144, 3).Value = -0.02418658464350926
$ws.Cells.Item(145, 3).Value = 0.2001520573564908
